# Backlog_6.xlsx edit
#
# The "Semana" (week) column (C) on both sheets was hard-coded from the
# shared string "Semana 06" to the plain number 6. Writing a literal
# number clears the t="s" shared-string reference on every affected
# cell; once the last reference to "Semana 06" disappears the shared
# string table itself shrinks (count 245->200, uniqueCount 31->30),
# which in turn shifts every other shared-string index above it (>14)
# down by one - exactly what the diff shows for columns B and I.

$wb  = $excel.ActiveWorkbook
$spn = $wb.Worksheets.Item(1)   # "SPN" sheet
$iti = $wb.Worksheets.Item(2)   # "ITI" sheet

# SPN: rows 2-14, column C ("Semana") -> literal 6
for ($r = 2; $r -le 14; $r++) {
    $spn.Range("C$r").Value = 6
}

# ITI: rows 2-33, column C ("Semana") -> literal 6
for ($r = 2; $r -le 33; $r++) {
    $iti.Range("C$r").Value = 6
}

# View/selection state: ITI loses the active tab, SPN gains it.
# Set ITI's own selection first (while it is still active) so that
# sheet's stored selection is updated, then activate SPN and select
# its new cell so SPN ends up as the active/selected tab.
$iti.Activate()
$iti.Range("A52:A54").Select()

$spn.Activate()
$spn.Range("C15").Select()
